$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 2).Value = 0.195583596214511
$ws.Cells.Item(2, 3).Value = 0.573080967402734
$ws.Cells.Item(2, 10).Value = 0.01787592008412198
$ws.Cells.Item(2, 16).Value = 0.1472134595162986
$ws.Cells.Item(2, 19).Value = 0.06624605678233439

# Row 3
$ws.Cells.Item(3, 2).Value = 0.00530035335689046
$ws.Cells.Item(3, 3).Value = 0.03003533568904593
$ws.Cells.Item(3, 10).Value = 0.03356890459363958
$ws.Cells.Item(3, 16).Value = 0.7544169611307421
$ws.Cells.Item(3, 19).Value = 0.176678445229682

# Row 5
$ws.Cells.Item(5, 10).Value = 0.125
$ws.Cells.Item(5, 16).Value = 0.625
$ws.Cells.Item(5, 19).Value = 0.25

# Row 6
$ws.Cells.Item(6, 2).Value = 0.05838323353293413
$ws.Cells.Item(6, 4).Value = 0.01197604790419162
$ws.Cells.Item(6, 5).Value = 0.002994011976047904
$ws.Cells.Item(6, 6).Value = 0.0658682634730539
$ws.Cells.Item(6, 10).Value = 0.2754491017964072
$ws.Cells.Item(6, 15).Value = 0.01047904191616766
$ws.Cells.Item(6, 17).Value = 0.1616766467065868
$ws.Cells.Item(6, 18).Value = 0.05389221556886228
$ws.Cells.Item(6, 19).Value = 0.3592814371257485

# Row 7
$ws.Cells.Item(7, 2).Value = 0.1282894736842105
$ws.Cells.Item(7, 4).Value = 0.01644736842105263
$ws.Cells.Item(7, 5).Value = 0.001644736842105263
$ws.Cells.Item(7, 6).Value = 0.04111842105263158
$ws.Cells.Item(7, 10).Value = 0.1299342105263158
$ws.Cells.Item(7, 15).Value = 0.0131578947368421
$ws.Cells.Item(7, 17).Value = 0.1825657894736842
$ws.Cells.Item(7, 18).Value = 0.08388157894736842
$ws.Cells.Item(7, 19).Value = 0.4029605263157895

# Row 8
$ws.Cells.Item(8, 2).Value = 0.1040118870728083
$ws.Cells.Item(8, 4).Value = 0.02451708766716196
$ws.Cells.Item(8, 5).Value = 0.002228826151560178
$ws.Cells.Item(8, 6).Value = 0.05943536404160475
$ws.Cells.Item(8, 10).Value = 0.09361069836552749
$ws.Cells.Item(8, 15).Value = 0.01783060921248143
$ws.Cells.Item(8, 17).Value = 0.175334323922734
$ws.Cells.Item(8, 18).Value = 0.08989598811292719
$ws.Cells.Item(8, 19).Value = 0.4331352154531947

# Row 9
$ws.Cells.Item(9, 2).Value = 0.1085141903171953
$ws.Cells.Item(9, 4).Value = 0.01335559265442404
$ws.Cells.Item(9, 6).Value = 0.05843071786310518
$ws.Cells.Item(9, 10).Value = 0.1035058430717863
$ws.Cells.Item(9, 15).Value = 0.02003338898163606
$ws.Cells.Item(9, 17).Value = 0.1886477462437396
$ws.Cells.Item(9, 18).Value = 0.08848080133555926
$ws.Cells.Item(9, 19).Value = 0.4190317195325542

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1106629694983615
$ws.Cells.Item(10, 4).Value = 0.02268716914544996
$ws.Cells.Item(10, 5).Value = 0.001008318628686665
$ws.Cells.Item(10, 6).Value = 0.0642803125787749
$ws.Cells.Item(10, 10).Value = 0.1041088984118982
$ws.Cells.Item(10, 15).Value = 0.0171414166876733
$ws.Cells.Item(10, 17).Value = 0.217544744139148
$ws.Cells.Item(10, 18).Value = 0.08721956138139653
$ws.Cells.Item(10, 19).Value = 0.375346609528611

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1428571428571428
$ws.Cells.Item(11, 10).Value = 0.07782515991471216
$ws.Cells.Item(11, 11).Value = 0.1908315565031983
$ws.Cells.Item(11, 12).Value = 0.5746268656716418
$ws.Cells.Item(11, 19).Value = 0.0138592750533049

# Row 12
$ws.Cells.Item(12, 7).Value = 0.7603603603603604
$ws.Cells.Item(12, 10).Value = 0.181981981981982
$ws.Cells.Item(12, 11).Value = 0.007207207207207207
$ws.Cells.Item(12, 12).Value = 0.03423423423423423
$ws.Cells.Item(12, 19).Value = 0.01621621621621622

# Row 15
$ws.Cells.Item(15, 6).Value = 0.01914580265095729
$ws.Cells.Item(15, 8).Value = 0.1222385861561119
$ws.Cells.Item(15, 9).Value = 0.08100147275405008
$ws.Cells.Item(15, 10).Value = 0.3681885125184094
$ws.Cells.Item(15, 11).Value = 0.06774668630338733
$ws.Cells.Item(15, 13).Value = 0.0117820324005891
$ws.Cells.Item(15, 14).Value = 0.001472754050073638
$ws.Cells.Item(15, 15).Value = 0.07658321060382917
$ws.Cells.Item(15, 19).Value = 0.251840942562592

# Row 16
$ws.Cells.Item(16, 6).Value = 0.02731411229135053
$ws.Cells.Item(16, 8).Value = 0.1714719271623672
$ws.Cells.Item(16, 9).Value = 0.08345978755690441
$ws.Cells.Item(16, 10).Value = 0.3975720789074355
$ws.Cells.Item(16, 11).Value = 0.1077389984825493
$ws.Cells.Item(16, 13).Value = 0.01669195751138088
$ws.Cells.Item(16, 15).Value = 0.0637329286798179
$ws.Cells.Item(16, 19).Value = 0.1320182094081942

# Row 17
$ws.Cells.Item(17, 6).Value = 0.02725366876310273
$ws.Cells.Item(17, 8).Value = 0.1586303284416492
$ws.Cells.Item(17, 9).Value = 0.08735150244584207
$ws.Cells.Item(17, 10).Value = 0.4409503843466108
$ws.Cells.Item(17, 11).Value = 0.09783368273934312
$ws.Cells.Item(17, 13).Value = 0.02026554856743536
$ws.Cells.Item(17, 15).Value = 0.05590496156533892
$ws.Cells.Item(17, 19).Value = 0.1118099231306778

# Row 18
$ws.Cells.Item(18, 6).Value = 0.01996672212978369
$ws.Cells.Item(18, 8).Value = 0.1331114808652246
$ws.Cells.Item(18, 9).Value = 0.09317803660565724
$ws.Cells.Item(18, 10).Value = 0.4442595673876872
$ws.Cells.Item(18, 11).Value = 0.09650582362728785
$ws.Cells.Item(18, 13).Value = 0.009983361064891847
$ws.Cells.Item(18, 14).Value = 0.001663893510815308
$ws.Cells.Item(18, 15).Value = 0.0632279534109817
$ws.Cells.Item(18, 19).Value = 0.1381031613976705

# Row 19
$ws.Cells.Item(19, 6).Value = 0.01541571897902451
$ws.Cells.Item(19, 8).Value = 0.2158200657063432
$ws.Cells.Item(19, 9).Value = 0.08011119535001264
$ws.Cells.Item(19, 10).Value = 0.3752843062926459
$ws.Cells.Item(19, 11).Value = 0.1079100328531716
$ws.Cells.Item(19, 13).Value = 0.01769016932019207
$ws.Cells.Item(19, 14).Value = 0.0005054334091483447
$ws.Cells.Item(19, 15).Value = 0.06823351023502654
$ws.Cells.Item(19, 19).Value = 0.1190295678544352
